$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.04247155220008958
$ws.Range("C2").Value = -0.470863874926198
$ws.Range("D2").Value = -1.596660741476967
$ws.Range("E2").Value = 1.095969070482927
$ws.Range("F2").Value = -0.517307604149827
$ws.Range("G2").Value = 0.5080998195400784
$ws.Range("H2").Value = -1.214309943970497
$ws.Range("I2").Value = 1.615420670929393
$ws.Range("J2").Value = -0.6191041902945924
$ws.Range("K2").Value = 0.7318558435072805
$ws.Range("B3").Value = -0.062631421304428
$ws.Range("C3").Value = -1.345037697590573
$ws.Range("D3").Value = 1.014063016659791
$ws.Range("E3").Value = -0.4422155547685414
$ws.Range("F3").Value = 0.6140069800538022
$ws.Range("G3").Value = -1.155893297811842
$ws.Range("H3").Value = 1.681859354759485
$ws.Range("I3").Value = -0.5444222316837559
$ws.Range("J3").Value = 0.8021691354082593
$ws.Range("K3").Value = 0.1802518522544978
$ws.Range("B4").Value = -1.450536909439388
$ws.Range("C4").Value = 0.8879577326944336
$ws.Range("D4").Value = -0.2739674861353097
$ws.Range("E4").Value = 0.6746710726576846
$ws.Range("F4").Value = -1.144760458437209
$ws.Range("G4").Value = 1.744722106086416
$ws.Range("H4").Value = -0.4870707038900127
$ws.Range("I4").Value = 0.8466743124126384
$ws.Range("J4").Value = 0.2313257232809164
$ws.Range("K4").Value = 0.7341089884009673
$ws.Range("B5").Value = 0.2835171128466949
$ws.Range("C5").Value = -0.3381708634900031
$ws.Range("D5").Value = 0.9662183807149013
$ws.Range("E5").Value = -1.182467334523142
$ws.Range("F5").Value = 1.736138086012954
$ws.Range("G5").Value = -0.4033296262926443
$ws.Range("H5").Value = 0.8815014573406833
$ws.Range("I5").Value = 0.2578901063332134
$ws.Range("J5").Value = 0.7796078291475662
$ws.Range("K5").Value = 0.1037568572541728
$ws.Range("B6").Value = -0.4258116923289144
$ws.Range("C6").Value = 0.9137991786852327
$ws.Range("D6").Value = -1.123820937538137
$ws.Range("E6").Value = 1.741067231350056
$ws.Range("F6").Value = -0.4144229093415396
$ws.Range("G6").Value = 0.8942681565657697
$ws.Range("H6").Value = 0.2666196653805725
$ws.Range("I6").Value = 0.782109773296418
$ws.Range("J6").Value = 0.1100327947258539
$ws.Range("K6").Value = 0.5363482025840406
$ws.Range("B7").Value = 0.9362271583182413
$ws.Range("C7").Value = -1.121217397975688
$ws.Range("D7").Value = 1.718999545696041
$ws.Range("E7").Value = -0.4204048896458946
$ws.Range("F7").Value = 0.8889249063833586
$ws.Range("G7").Value = 0.2562511320102062
$ws.Range("H7").Value = 0.7736324511727942
$ws.Range("I7").Value = 0.1022604864190431
$ws.Range("J7").Value = 0.527724229051072
$ws.Range("K7").Value = 0.5878314504560218
$ws.Range("B8").Value = -1.099792826518468
$ws.Range("C8").Value = 1.839177394495253
$ws.Range("D8").Value = -0.5095992340596777
$ws.Range("E8").Value = 0.8591979025347809
$ws.Range("F8").Value = 0.2648772520401163
$ws.Range("G8").Value = 0.747692277747359
$ws.Range("H8").Value = 0.07948192339081561
$ws.Range("I8").Value = 0.5141061939024542
$ws.Range("J8").Value = 0.569310653270846
$ws.Range("K8").Value = 0.6926620895998143
$ws.Range("B9").Value = 1.796856762174863
$ws.Range("C9").Value = -0.5388037079154109
$ws.Range("D9").Value = 0.8561773027298597
$ws.Range("E9").Value = 0.2460159684286359
$ws.Range("F9").Value = 0.7273611466821548
$ws.Range("G9").Value = 0.06495176795656782
$ws.Range("H9").Value = 0.4974651484101862
$ws.Range("I9").Value = 0.5517268078420958
$ws.Range("J9").Value = 0.6761882060697715
$ws.Range("K9").Value = -0.2163793123768544
$ws.Range("B10").Value = -0.1986059672975008
$ws.Range("C10").Value = 0.9315483645137967
$ws.Range("D10").Value = 0.05492268211563681
$ws.Range("E10").Value = 0.7520550592688879
$ws.Range("F10").Value = 0.07433117151286378
$ws.Range("G10").Value = 0.4426543497987951
$ws.Range("H10").Value = 0.532832716953426
$ws.Range("I10").Value = 0.6613242052547258
$ws.Range("J10").Value = -0.2450910230062252
$ws.Range("K10").Value = 0.5297973106668776
$ws.Range("B11").Value = 1.382243236504047
$ws.Range("C11").Value = 0.1020333843279952
$ws.Range("D11").Value = 0.5117718284780797
$ws.Range("E11").Value = 0.1102762525787266
$ws.Range("F11").Value = 0.4406035923592023
$ws.Range("G11").Value = 0.4562644775173959
$ws.Range("H11").Value = 0.6339248108794424
$ws.Range("I11").Value = -0.2713800732039505
$ws.Range("J11").Value = 0.4864000199535451
$ws.Range("K11").Value = 0.2305062539156956
$ws.Range("B12").Value = 0.4143589712515336
$ws.Range("C12").Value = 0.645060568213604
$ws.Range("D12").Value = -0.06880216923746499
$ws.Range("E12").Value = 0.4687978311297695
$ws.Range("F12").Value = 0.4922897887097893
$ws.Range("G12").Value = 0.5963832406752287
$ws.Range("H12").Value = -0.2769224122913118
$ws.Range("I12").Value = 0.4905753488922937
$ws.Range("J12").Value = 0.2193215401759246
$ws.Range("B13").Value = 0.8806510593214452
$ws.Range("C13").Value = 0.01587657163684458
$ws.Range("D13").Value = 0.3231974644960524
$ws.Range("E13").Value = 0.5022435315690319
$ws.Range("F13").Value = 0.6117484103747737
$ws.Range("G13").Value = -0.3157707472487348
$ws.Range("H13").Value = 0.4747295398651078
$ws.Range("I13").Value = 0.2109873117084238
$ws.Range("B14").Value = 0.3269007581182274
$ws.Range("C14").Value = 0.4402862389072476
$ws.Range("D14").Value = 0.3397793901194651
$ws.Range("E14").Value = 0.6407073608399754
$ws.Range("F14").Value = -0.2800278747938312
$ws.Range("G14").Value = 0.4434018934814807
$ws.Range("H14").Value = 0.2077622620068982
$ws.Range("B15").Value = 0.687111416311362
$ws.Range("C15").Value = 0.3574352576573546
$ws.Range("D15").Value = 0.5460851217289808
$ws.Range("E15").Value = -0.2478104865453511
$ws.Range("F15").Value = 0.455790019309198
$ws.Range("G15").Value = 0.1868984584576193
$ws.Range("B16").Value = 0.5971116170735665
$ws.Range("C16").Value = 0.6311093651981943
$ws.Range("D16").Value = -0.3642950649033653
$ws.Range("E16").Value = 0.4708952909610505
$ws.Range("F16").Value = 0.2101374940836094
$ws.Range("B17").Value = 0.7959694422322816
$ws.Range("C17").Value = -0.350757792655449
$ws.Range("D17").Value = 0.4079223719710875
$ws.Range("E17").Value = 0.2201756597651073
$ws.Range("B18").Value = -0.09690875079004102
$ws.Range("C18").Value = 0.5072404221531239
$ws.Range("D18").Value = 0.1085991175498651
$ws.Range("B19").Value = 0.5513001133925729
$ws.Range("C19").Value = 0.130019622424466
$ws.Range("B20").Value = 0.3662627537369125
